# "Putting art asset list" — populate the Art List sheet with the list of
# art assets, make that sheet the active tab, and leave the selection on H14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Art List")

$assets = @(
    'Trex idle',
    'Pterodactyl/ animation',
    '13 Cactuses',
    'Cloud',
    'Ground',
    'Reset Button',
    'Trex Logo',
    'Trex run animation',
    'Trex death frame ',
    'Trex crouch/run animation',
    'Chicken run ',
    'Chicken idle',
    'Chicken jump ',
    'chicken coup',
    'chicken nest',
    'egg',
    'seed ',
    'background',
    'ground',
    'farmer',
    'dog enemy animation',
    'sign ',
    'chicken death'
)

for ($i = 0; $i -lt $assets.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $assets[$i]
}

$ws.Activate()
$ws.Range("H14").Select() | Out-Null
